# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 7;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 8;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 9;  DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 11; DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 13; DAMSLTag = "qy"; DialogAct = "Yes-No-Question" },
    @{ Row = 14; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 17; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 22; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 23; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 39; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 41; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 45; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.DAMSLTag
    $ws.Range("J" + $u.Row).Value = $u.DialogAct
}
